$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 67.54763293266296, 114.03125),
    @(39.78, 0.16, 71.40265464782715, 125.984375),
    @(39.78, 0.16, 63.66919302940369, 110.875),
    @(39.78, 0.16, 74.6501157283783, 126.1875),
    @(50.43, 0.16, 80.41346979141235, 127.59375),
    @(50.43, 0.16, 86.68496823310852, 144.9375),
    @(99.67, 0.16, 85.03541254997253, 141.15625),
    @(10.65, 0.16, 67.56848335266113, 118.265625),
    @(39.78, 0.16, 63.21611595153809, 100.875),
    @(20.78, 0.16, 66.74141430854797, 113.921875)
)

$startRow = 96
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
